$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.245.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.555.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.37%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.24%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.58%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.771.95"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.552.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.22%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.32%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.249.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.57%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "185.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.31%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.49%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.59"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.46%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.59%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.84%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.03"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.18%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.75%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.83%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.084.92"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0148"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.493"

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.765"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -10.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.75"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.75%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.686.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.37%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.22%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.80%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.15%  "
